# Data/data_agency.xlsx - "VV function - var in %?" revision
#
# The variance parameters for three "tnorm_0_1" variables were switched to
# "posnorm" distributions expressed as percentages (0-100) instead of
# fractions (0-1), the descriptions were renamed from "Variance in ..." to
# "Coefficient of variation ...", and a couple of unit labels were tidied up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: labor_wage_Eur_per_h_brutto -> unit label fix ---
$ws.Range("F5").Value = "Euro"

# --- Row 6: tree_establishment_material_costs_per_tree -> unit label fix ---
$ws.Range("F6").Value = "Eur/tree"

# --- Row 8: tree_labor_establishment_var_per_tree ---
# distribution tnorm_0_1 -> posnorm, values scaled from fraction to percent,
# unit absolute -> percent, description renamed to "Coefficient of variation"
$ws.Range("B8").Value = "posnorm"
$ws.Range("C8").Value = 10
$ws.Range("E8").Value = 20
$ws.Range("F8").Value = "percent"
$ws.Range("G8").Value = "Coefficient of variation for 1rst tree planting"

# --- Row 10: tree_mainteance_costs_fertiliser_mean_per_tree -> unit label fix ---
$ws.Range("F10").Value = "Eur/tree"

# --- Row 11: tree_mainteance_costs_fertiliser_var ---
$ws.Range("B11").Value = "posnorm"
$ws.Range("C11").Value = 5
$ws.Range("E11").Value = 30
$ws.Range("F11").Value = "percent"
$ws.Range("G11").Value = "Coefficient of variation in fertiliser costs"

# --- Row 22: labor_fruit_replanting_var_per_tree ---
$ws.Range("B22").Value = "posnorm"
$ws.Range("C22").Value = 10
$ws.Range("E22").Value = 20
$ws.Range("G22").Value = "Coefficient of variation for labor of 1rst tree planting"

# G2 loses its extra fill/border formatting (reverts to the default style)
$ws.Range("G2").ClearFormats()

# Selection moved from A36 to A31
[void]$ws.Range("A31").Select()
